$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J5 previously referenced a shared string ("?"); now a numeric value
$ws.Range("J5").Value = 16430

# Fill in previously empty J-column cells with corrected postal values
$ws.Range("J39").Value = 41271
$ws.Range("J40").Value = 40294
$ws.Range("J41").Value = 41181
$ws.Range("J44").Value = 17530
$ws.Range("J45").Value = 15720
$ws.Range("J47").Value = 17530
$ws.Range("J48").Value = 45154
$ws.Range("J49").Value = 16670
$ws.Range("J50").Value = 17530
$ws.Range("J51").Value = 15820
$ws.Range("J53").Value = 17530
$ws.Range("J54").Value = 16340
$ws.Range("J56").Value = 16370
$ws.Range("J57").Value = 40554
$ws.Range("J58").Value = 15710
$ws.Range("J60").Value = 16640
$ws.Range("J63").Value = 16810
$ws.Range("J64").Value = 43284
$ws.Range("J65").Value = 78821
$ws.Range("J66").Value = 16120
$ws.Range("J68").Value = 16710
$ws.Range("J69").Value = 16770
$ws.Range("J71").Value = 40562
$ws.Range("J72").Value = 16830
$ws.Range("J74").Value = 40288
$ws.Range("J75").Value = 16740
$ws.Range("J77").Value = 15720
$ws.Range("J78").Value = 16750
$ws.Range("J79").Value = 45161
$ws.Range("J81").Value = 16922
$ws.Range("J82").Value = 40195
$ws.Range("J83").Value = 17132
$ws.Range("J84").Value = 16451
$ws.Range("J86").Value = 40973
$ws.Range("J88").Value = 16810
$ws.Range("J89").Value = 16810
$ws.Range("J90").Value = 16810
$ws.Range("J91").Value = 15710
$ws.Range("J92").Value = 16660
$ws.Range("J93").Value = 40561
$ws.Range("J94").Value = 44184
$ws.Range("J95").Value = 15560
$ws.Range("J96").Value = 17630
$ws.Range("J97").Value = 46196
$ws.Range("J99").Value = 40294
$ws.Range("J101").Value = 41374
$ws.Range("J102").Value = 45154
$ws.Range("J103").Value = 15610
$ws.Range("J105").Value = 41181
$ws.Range("J107").Value = 45363
$ws.Range("J108").Value = 40562
$ws.Range("J109").Value = 44151
$ws.Range("J111").Value = 60183
$ws.Range("J112").Value = 41374
$ws.Range("J115").Value = 41362
$ws.Range("J116").Value = 17530
$ws.Range("J117").Value = 15530
$ws.Range("J121").Value = 17630
$ws.Range("J123").Value = 17620
$ws.Range("J127").Value = 40294
$ws.Range("J128").Value = 44151
$ws.Range("J129").Value = 15720
$ws.Range("J130").Value = 51137
$ws.Range("J131").Value = 12620
$ws.Range("J132").Value = 61372
$ws.Range("J133").Value = 57731
$ws.Range("J134").Value = 67225
$ws.Range("J135").Value = 61265
$ws.Range("J136").Value = 61265
$ws.Range("J137").Value = 64411
$ws.Range("J138").Value = 57731
$ws.Range("J139").Value = 46125
$ws.Range("J140").Value = 61352
$ws.Range("J141").Value = 57731
$ws.Range("J142").Value = 57731
$ws.Range("J144").Value = 64411
$ws.Range("J145").Value = 64411
$ws.Range("J146").Value = 41282
$ws.Range("J147").Value = 16810
$ws.Range("J148").Value = 61328
$ws.Range("J149").Value = 57731
$ws.Range("J151").Value = 16810
$ws.Range("J156").Value = 57731
$ws.Range("J157").Value = 57731
$ws.Range("J166").Value = 16810
$ws.Range("J170").Value = 15560
$ws.Range("J175").Value = 15720
$ws.Range("J201").Value = 15560
$ws.Range("J216").Value = 41374
$ws.Range("J221").Value = 40294
$ws.Range("J222").Value = 40294
$ws.Range("J225").Value = 15560
$ws.Range("J234").Value = 41374
$ws.Range("J235").Value = 40294
$ws.Range("J237").Value = 57731
$ws.Range("J244").Value = 41374
$ws.Range("J254").Value = 16810
$ws.Range("J262").Value = 45154
$ws.Range("J266").Value = 16810
$ws.Range("J268").Value = 16810
$ws.Range("J271").Value = 40294
$ws.Range("J273").Value = 45154
$ws.Range("J277").Value = 16810
$ws.Range("J279").Value = 41374
$ws.Range("J281").Value = 51137
$ws.Range("J282").Value = 41181
$ws.Range("J285").Value = 16810
$ws.Range("J289").Value = 15610
$ws.Range("J294").Value = 16810

# Update the active sheet view/selection to match the saved state
# (scrolled so row 91 / column C is the top-left visible cell, with I95 selected)
$ws.Application.ActiveWindow.ScrollRow = 91
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("I95").Select()
